# 財產申報表 (property declaration) — #5: property boat&car done
# Bring the 汽車 (car) sheet up to parity with the other sheets: give it a
# real header row and append the legislator/source-file metadata columns
# (H:N) that every other sheet already carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: column headers (bold/bordered style already applied) ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2: data (existing car record, now with metadata columns) ---
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2012-02-29"
$ws.Cells.Item(2, 11).Value = "林岱樺"
$ws.Cells.Item(2, 12).Value = 904
$ws.Cells.Item(2, 13).Value = "tmp3bff1"
$ws.Cells.Item(2, 14).Value = 29
